$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.070371838301993
$ws.Cells.Item(2, 4).Value = 1.068807906730887
$ws.Cells.Item(2, 5).Value = 1.074087081601978
$ws.Cells.Item(2, 6).Value = 1.078532240469478
$ws.Cells.Item(2, 9).Value = 1.051147693843331
$ws.Cells.Item(2, 10).Value = 1.075301829498439
$ws.Cells.Item(2, 11).Value = 1.071512099279177
$ws.Cells.Item(2, 12).Value = 1.076777223258366
$ws.Cells.Item(2, 13).Value = 1.081210666265848
$ws.Cells.Item(2, 14).Value = 1.076828881435171

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.072427629661114
$ws.Cells.Item(3, 4).Value = 1.070403482192971
$ws.Cells.Item(3, 5).Value = 1.07609627420989
$ws.Cells.Item(3, 6).Value = 1.08027992115699
$ws.Cells.Item(3, 9).Value = 1.051771661125545
$ws.Cells.Item(3, 10).Value = 1.077009949197536
$ws.Cells.Item(3, 11).Value = 1.072921789672596
$ws.Cells.Item(3, 12).Value = 1.07860054855899
$ws.Cells.Item(3, 13).Value = 1.082773983907301
$ws.Cells.Item(3, 14).Value = 1.078539426860164

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.073751960065182
$ws.Cells.Item(4, 4).Value = 1.071430334281427
$ws.Cells.Item(4, 5).Value = 1.077390822473071
$ws.Cells.Item(4, 6).Value = 1.081405474269382
$ws.Cells.Item(4, 9).Value = 1.052171205811199
$ws.Cells.Item(4, 10).Value = 1.078109118681496
$ws.Cells.Item(4, 11).Value = 1.07382788210049
$ws.Cells.Item(4, 12).Value = 1.079774431592579
$ws.Cells.Item(4, 13).Value = 1.083779787891898
$ws.Cells.Item(4, 14).Value = 1.079640157290868

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.074307325745585
$ws.Cells.Item(5, 4).Value = 1.071860707818625
$ws.Cells.Item(5, 5).Value = 1.077933754843347
$ws.Cells.Item(5, 6).Value = 1.081877409704638
$ws.Cells.Item(5, 9).Value = 1.052338179086276
$ws.Cells.Item(5, 10).Value = 1.078569775849787
$ws.Cells.Item(5, 11).Value = 1.074207372282963
$ws.Cells.Item(5, 12).Value = 1.080266539659569
$ws.Cells.Item(5, 13).Value = 1.084201270829418
$ws.Cells.Item(5, 14).Value = 1.080101468645176

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.074400493882425
$ws.Cells.Item(6, 4).Value = 1.071932892995021
$ws.Cells.Item(6, 5).Value = 1.078024840435209
$ws.Cells.Item(6, 6).Value = 1.081956577329541
$ws.Cells.Item(6, 9).Value = 1.052366156588952
$ws.Cells.Item(6, 10).Value = 1.078647038901219
$ws.Cells.Item(6, 11).Value = 1.074271007084214
$ws.Cells.Item(6, 12).Value = 1.080349085911577
$ws.Cells.Item(6, 13).Value = 1.084271960772041
$ws.Cells.Item(6, 14).Value = 1.080178841419003

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.07375938628096
$ws.Cells.Item(7, 4).Value = 1.07143609008924
$ws.Cells.Item(7, 5).Value = 1.077398082214183
$ws.Cells.Item(7, 6).Value = 1.081411785160482
$ws.Cells.Item(7, 9).Value = 1.052173440812196
$ws.Cells.Item(7, 10).Value = 1.07811527960323
$ws.Cells.Item(7, 11).Value = 1.073832958459797
$ws.Cells.Item(7, 12).Value = 1.079781012598852
$ws.Cells.Item(7, 13).Value = 1.083785425068092
$ws.Cells.Item(7, 14).Value = 1.079646326961818

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.071067843651303
$ws.Cells.Item(8, 4).Value = 1.069348311100486
$ws.Cells.Item(8, 5).Value = 1.074767261050243
$ws.Cells.Item(8, 6).Value = 1.079123992546818
$ws.Cells.Item(8, 9).Value = 1.051359443678974
$ws.Cells.Item(8, 10).Value = 1.075880375114769
$ws.Cells.Item(8, 11).Value = 1.07198978150862
$ws.Cells.Item(8, 12).Value = 1.077394669092008
$ws.Cells.Item(8, 13).Value = 1.081740205187956
$ws.Cells.Item(8, 14).Value = 1.077408248652636

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.06627839578289
$ws.Cells.Item(9, 4).Value = 1.065625478636546
$ws.Cells.Item(9, 5).Value = 1.070087719576967
$ws.Cells.Item(9, 6).Value = 1.075050789085273
$ws.Cells.Item(9, 9).Value = 1.049892376193196
$ws.Cells.Item(9, 10).Value = 1.071894293963722
$ws.Cells.Item(9, 11).Value = 1.06869434929565
$ws.Cells.Item(9, 12).Value = 1.073142966564688
$ws.Cells.Item(9, 13).Value = 1.078091028131786
$ws.Cells.Item(9, 14).Value = 1.073416506809144

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.063052137988969
$ws.Cells.Item(10, 4).Value = 1.063112562418062
$ws.Cells.Item(10, 5).Value = 1.066936787226883
$ws.Cells.Item(10, 6).Value = 1.072305610148356
$ws.Cells.Item(10, 9).Value = 1.048891633047746
$ws.Cells.Item(10, 10).Value = 1.069203027407591
$ws.Cells.Item(10, 11).Value = 1.066464007953935
$ws.Cells.Item(10, 12).Value = 1.070275375572439
$ws.Cells.Item(10, 13).Value = 1.075626318126554
$ws.Cells.Item(10, 14).Value = 1.070721418345811

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.061646795725028
$ws.Cells.Item(11, 4).Value = 1.062016741853732
$ws.Cells.Item(11, 5).Value = 1.06556457140065
$ws.Cells.Item(11, 6).Value = 1.071109511996912
$ws.Cells.Item(11, 9).Value = 1.04845275584257
$ws.Cells.Item(11, 10).Value = 1.068029265797396
$ws.Cells.Item(11, 11).Value = 1.065490002878363
$ws.Cells.Item(11, 12).Value = 1.069025426618898
$ws.Cells.Item(11, 13).Value = 1.074551159194673
$ws.Cells.Item(11, 14).Value = 1.069545989859498

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.061123496083204
$ws.Cells.Item(12, 4).Value = 1.061608517147413
$ws.Cells.Item(12, 5).Value = 1.065053654957015
$ws.Cells.Item(12, 6).Value = 1.07066408231173
$ws.Cells.Item(12, 9).Value = 1.048288889648975
$ws.Cells.Item(12, 10).Value = 1.067591979806066
$ws.Cells.Item(12, 11).Value = 1.065126946002531
$ws.Cells.Item(12, 12).Value = 1.068559864228658
$ws.Cells.Item(12, 13).Value = 1.074150577346193
$ws.Cells.Item(12, 14).Value = 1.069108082871902

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.061235804762748
$ws.Cells.Item(13, 4).Value = 1.061696136980427
$ws.Cells.Item(13, 5).Value = 1.065163303811834
$ws.Cells.Item(13, 6).Value = 1.070759680837144
$ws.Cells.Item(13, 9).Value = 1.04832407808969
$ws.Cells.Item(13, 10).Value = 1.067685838436351
$ws.Cells.Item(13, 11).Value = 1.065204880773007
$ws.Cells.Item(13, 12).Value = 1.068659787225837
$ws.Cells.Item(13, 13).Value = 1.074236559198048
$ws.Cells.Item(13, 14).Value = 1.069202074792208

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.061603566212928
$ws.Cells.Item(14, 4).Value = 1.061983022283477
$ws.Cells.Item(14, 5).Value = 1.065522363873548
$ws.Cells.Item(14, 6).Value = 1.071072716191864
$ws.Cells.Item(14, 9).Value = 1.048439227985833
$ws.Cells.Item(14, 10).Value = 1.067993146285423
$ws.Cells.Item(14, 11).Value = 1.065460018515785
$ws.Cells.Item(14, 12).Value = 1.068986969310582
$ws.Cells.Item(14, 13).Value = 1.074518072033672
$ws.Cells.Item(14, 14).Value = 1.069509819053674

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.061829983695438
$ws.Cells.Item(15, 4).Value = 1.062159623449442
$ws.Cells.Item(15, 5).Value = 1.06574343061256
$ws.Cells.Item(15, 6).Value = 1.071265434863611
$ws.Cells.Item(15, 9).Value = 1.04851006289078
$ws.Cells.Item(15, 10).Value = 1.068182315715255
$ws.Cells.Item(15, 11).Value = 1.065617048423096
$ws.Cells.Item(15, 12).Value = 1.069188386914426
$ws.Cells.Item(15, 13).Value = 1.074691358833973
$ws.Cells.Item(15, 14).Value = 1.069699257125795

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.063145226722965
$ws.Cells.Item(16, 4).Value = 1.063185123435424
$ws.Cells.Item(16, 5).Value = 1.06702768837866
$ws.Cells.Item(16, 6).Value = 1.072384832315719
$ws.Cells.Item(16, 9).Value = 1.048920641784984
$ws.Cells.Item(16, 10).Value = 1.069280745795582
$ws.Cells.Item(16, 11).Value = 1.06652847319787
$ws.Cells.Item(16, 12).Value = 1.070358153563054
$ws.Cells.Item(16, 13).Value = 1.075697503437981
$ws.Cells.Item(16, 14).Value = 1.070799247102828

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.063967981232357
$ws.Cells.Item(17, 4).Value = 1.063826307206335
$ws.Cells.Item(17, 5).Value = 1.067831144558291
$ws.Cells.Item(17, 6).Value = 1.073084993037173
$ws.Cells.Item(17, 9).Value = 1.049176691688902
$ws.Cells.Item(17, 10).Value = 1.069967483231665
$ws.Cells.Item(17, 11).Value = 1.067097956703534
$ws.Cells.Item(17, 12).Value = 1.071089681007281
$ws.Cells.Item(17, 13).Value = 1.076326489375617
$ws.Cells.Item(17, 14).Value = 1.07148695978484

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.064447076850333
$ws.Cells.Item(18, 4).Value = 1.064199556854231
$ws.Cells.Item(18, 5).Value = 1.068299032864337
$ws.Cells.Item(18, 6).Value = 1.073492671386086
$ws.Cells.Item(18, 9).Value = 1.049325506844281
$ws.Cells.Item(18, 10).Value = 1.070367234874085
$ws.Cells.Item(18, 11).Value = 1.067429332472708
$ws.Cells.Item(18, 12).Value = 1.071515573659798
$ws.Cells.Item(18, 13).Value = 1.076692603777425
$ws.Cells.Item(18, 14).Value = 1.071887279120444

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.064610300992031
$ws.Cells.Item(19, 4).Value = 1.064326700371804
$ws.Cells.Item(19, 5).Value = 1.068458443868445
$ws.Cells.Item(19, 6).Value = 1.073631559071851
$ws.Cells.Item(19, 9).Value = 1.049376158739619
$ws.Cells.Item(19, 10).Value = 1.070503403332793
$ws.Cells.Item(19, 11).Value = 1.067542189160418
$ws.Cells.Item(19, 12).Value = 1.071660658242297
$ws.Cells.Item(19, 13).Value = 1.076817310889079
$ws.Cells.Item(19, 14).Value = 1.072023640953983

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.063879790810325
$ws.Cells.Item(20, 4).Value = 1.06375759121709
$ws.Cells.Item(20, 5).Value = 1.067745019574184
$ws.Cells.Item(20, 6).Value = 1.073009946448338
$ws.Cells.Item(20, 9).Value = 1.049149275335871
$ws.Cells.Item(20, 10).Value = 1.069893886893041
$ws.Cells.Item(20, 11).Value = 1.067036938876128
$ws.Cells.Item(20, 12).Value = 1.071011277505088
$ws.Cells.Item(20, 13).Value = 1.076259084213939
$ws.Cells.Item(20, 14).Value = 1.071413258930973

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.061495305634381
$ws.Cells.Item(21, 4).Value = 1.061898574751617
$ws.Cells.Item(21, 5).Value = 1.065416663452532
$ws.Cells.Item(21, 6).Value = 1.070980566958617
$ws.Cells.Item(21, 9).Value = 1.048405342724748
$ws.Cells.Item(21, 10).Value = 1.067902687899479
$ws.Cells.Item(21, 11).Value = 1.065384922041996
$ws.Cells.Item(21, 12).Value = 1.068890657834875
$ws.Cells.Item(21, 13).Value = 1.07443520743675
$ws.Cells.Item(21, 14).Value = 1.069419232206447

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.05998858510881
$ws.Cells.Item(22, 4).Value = 1.060722847949998
$ws.Cells.Item(22, 5).Value = 1.063945689024777
$ws.Cells.Item(22, 6).Value = 1.069697971744069
$ws.Cells.Item(22, 9).Value = 1.047932691350431
$ws.Cells.Item(22, 10).Value = 1.06664321168636
$ws.Cells.Item(22, 11).Value = 1.06433888278511
$ws.Cells.Item(22, 12).Value = 1.067549941902396
$ws.Cells.Item(22, 13).Value = 1.073281390437941
$ws.Cells.Item(22, 14).Value = 1.068157967392642

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.060788050032067
$ws.Cells.Item(23, 4).Value = 1.061346786032599
$ws.Cells.Item(23, 5).Value = 1.064726160429722
$ws.Cells.Item(23, 6).Value = 1.070378540038115
$ws.Cells.Item(23, 9).Value = 1.048183723075929
$ws.Cells.Item(23, 10).Value = 1.067311608764144
$ws.Cells.Item(23, 11).Value = 1.064894114282961
$ws.Cells.Item(23, 12).Value = 1.068261393639023
$ws.Cells.Item(23, 13).Value = 1.07389373113416
$ws.Cells.Item(23, 14).Value = 1.068827313670941

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.063919642746193
$ws.Cells.Item(24, 4).Value = 1.063788643306417
$ws.Cells.Item(24, 5).Value = 1.067783938076018
$ws.Cells.Item(24, 6).Value = 1.073043858975191
$ws.Cells.Item(24, 9).Value = 1.049161665256487
$ws.Cells.Item(24, 10).Value = 1.069927144414262
$ws.Cells.Item(24, 11).Value = 1.067064512659196
$ws.Cells.Item(24, 12).Value = 1.07104670712787
$ws.Cells.Item(24, 13).Value = 1.076289544059914
$ws.Cells.Item(24, 14).Value = 1.071446563681689

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.067522305442571
$ws.Cells.Item(25, 4).Value = 1.066593276513492
$ws.Cells.Item(25, 5).Value = 1.071302861086399
$ws.Cells.Item(25, 6).Value = 1.07610892805995
$ws.Cells.Item(25, 9).Value = 1.050275596042892
$ws.Cells.Item(25, 10).Value = 1.072930638201145
$ws.Cells.Item(25, 11).Value = 1.069552075263563
$ws.Cells.Item(25, 12).Value = 1.079039942711277
$ws.Cells.Item(25, 13).Value = 1.074454322774254
$ws.Cells.Item(25, 14).Value = 1.074454322774254
